$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Insert Ambulance Map" section ---------------------------------------
# Row 6: section header, highlighted with the accent fill
$ws.Range("A6").Value = "Insert Ambulance Map"
$ws.Range("A6:B6").Interior.ThemeColor = 10

# Rows 7-9: response code / message pairs
$ws.Range("A7").Value = 0
$ws.Range("B7").Value = "Insertion Successful"

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Failed To insert because the car already has resources assigned"

$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "Failed to insert because car is already assigned elsewhere"

# --- "Delete Ambulance Map" section ---------------------------------------
# Row 10: section header, highlighted with the accent fill
$ws.Range("A10").Value = "Delete Ambulance Map"
$ws.Range("A10:B10").Interior.ThemeColor = 10

# Rows 11-12: response code / message pairs
$ws.Range("A11").Value = 0
$ws.Range("B11").Value = "Deleted Successfully"

$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Failed to delete"

# Match the saved selection state from the authored workbook
$ws.Range("B13").Select()
